# Credit.xlsx — add the "SQLite in Unity" reference row to the Referensi
# sheet, widen the Sumber Tautan (link) column to fit it, and leave the
# Referensi tab as the active/selected sheet (matching the author's final
# view state).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Referensi"
$ws2 = $wb.Worksheets.Item(2)   # "Aset"

# --- New row 3 on "Referensi": Adorable73 / SQLite-in-Unity tutorial ---
$ws1.Range("A3").Value = "Adorable73"
$ws1.Range("B3").Value = "https://answers.unity.com/questions/743400/database-sqlite-setup-for-unity.html"
$ws1.Range("C3").Value = "SQLite in Unity"
$ws1.Range("D3").Value = "Setup SQLite in Unity tutorial"

# --- Widen column B so the long URL is readable ---
$ws1.Columns.Item(2).ColumnWidth = 75.77734375

# --- Restore sheet2's remembered selection before switching away from it ---
$null = $ws2.Activate()
$null = $ws2.Range("A3").Select()

# --- Make "Referensi" the active tab again, selecting the newly added cell ---
$null = $ws1.Activate()
$null = $ws1.Range("D3").Select()
